# Append: 2025-10-20 18:26 JST
# Two new job postings were scraped; insert them into the "ランサーズ" sheet
# at their sorted position (by score), and refresh the "取得日時" timestamp
# on every existing row to the new scrape time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-20 18:26:29"

# --- 1) Insert the two new rows at their final positions -------------------
# Row 21: a new posting (FPC adapter) slots in just above the old row 21.
$ws.Rows("21:21").Insert()
# Row 24 (after the previous insert shifted everything below +1): a second
# new posting (the Chinese-language escort ad) slots in just above what is
# now row 24 (old row 23).
$ws.Rows("24:24").Insert()

# --- 2) Populate the newly inserted rows ------------------------------------
$ws.Cells.Item(21, 1).Value = $newTimestamp
$ws.Cells.Item(21, 2).Value = "【急募】既製品へのファームウェア書き込み用FPCアダプタ製作依頼"
$ws.Cells.Item(21, 3).Value = "システム開発"
$ws.Cells.Item(21, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(21, 5).Value = "期限情報なし"
$ws.Cells.Item(21, 6).Value = "https://www.lancers.jp/work/detail/5416679"
$ws.Cells.Item(21, 7).Value = 18

$ws.Cells.Item(24, 1).Value = $newTimestamp
$ws.Cells.Item(24, 2).Value = "台灣高檔外送茶Gleezy搜索賬號id3p6688高檔外約/網美/寫真模特兒/百萬粉Coser怎麼約"
$ws.Cells.Item(24, 3).Value = "システム開発"
$ws.Cells.Item(24, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(24, 5).Value = "期限情報なし"
$ws.Cells.Item(24, 6).Value = "https://www.lancers.jp/work/detail/5417087"
$ws.Cells.Item(24, 7).Value = 10

# --- 3) Refresh the scrape timestamp on every row (2..26) -------------------
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 4) Rebuild the hyperlinks on column F in row order ---------------------
# Row inserts do not renumber existing hyperlink anchors in this engine, so
# drop every hyperlink on the sheet and recreate them fresh, in row order,
# against the (already correct) URL text now sitting in column F.
$ws.Range("A1").Hyperlinks.Delete()

for ($r = 2; $r -le 26; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value()
    $ws.Hyperlinks.Add($cell, $url)
    $cell.Style = "Hyperlink"
}
